$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.339.51"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.90"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.15"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6305"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9984"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07579"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2916"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.846.38"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.004"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6784"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001047"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.18"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.324.00"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.38"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9981"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.448"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9977"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.32"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1394"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.434"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.63"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.425"
$ws.Range("E28").Value = "  +5.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.473"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05672"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.042"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.823"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7000"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.573"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.242.09"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.716"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.405"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9023"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9978"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.006.68"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.53"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.127"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1164"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.984"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3951"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.672"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000113"
$ws.Range("E51").Value = "  -4.99%  "
